$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per repull of data / mean calculation fix
$updates = @{
    2  = -5
    4  = 0
    9  = -1
    15 = -1
    19 = -1
    21 = -2
    22 = 1
    23 = 0
    27 = 2
    31 = 2
    37 = 1
    41 = 2
    46 = -3
    58 = 1
    61 = -2
    66 = -4
    68 = -3
    69 = 0
    70 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
